$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the text values currently held in A4 and B4 (hyperlinks stay attached
# to the cells as-is; only the displayed/stored text is swapped).
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2

$ws.Range("A4").Value = $b4
$ws.Range("B4").Value = $a4

# Move the active selection to B5, matching the saved cursor position.
$ws.Range("B5").Select()

$wb.Save()
